$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

$ws.Range("D3").Value = "70% (Chưa có cn Cập Nhật)"
$ws.Range("D4").Value = "70% (Chưa có cn Cập Nhật)"
$ws.Range("D5").Value = "70% (Chưa có cn Cập Nhật)"
$ws.Range("E6").Value = "100% (9/06/2010)"
